$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 4.914230582351781
$ws.Range("D2").Value = 7.39804222745583
$ws.Range("E2").Value = 12.28299269679532
$ws.Range("F2").Value = 39.33569652803837
$ws.Range("G2").Value = 46.32005703397751
$ws.Range("H2").Value = 18.45746421273209
$ws.Range("J2").Value = 9.868275483724016
$ws.Range("K2").Value = 17.0168400435231
$ws.Range("M2").Value = 18.63935577571864
$ws.Range("N2").Value = 19.19277078021837
# Row 3
$ws.Range("C3").Value = 4.752483729100021
$ws.Range("D3").Value = 7.393720686979336
$ws.Range("E3").Value = 12.30034240665037
$ws.Range("F3").Value = 39.32018844834531
$ws.Range("G3").Value = 46.19323033242995
$ws.Range("H3").Value = 18.49423114194269
$ws.Range("J3").Value = 9.892613564710095
$ws.Range("K3").Value = 16.65498122017688
$ws.Range("M3").Value = 18.5109125860743
$ws.Range("N3").Value = 19.26531851198388
# Row 4
$ws.Range("C4").Value = 4.651922307078825
$ws.Range("D4").Value = 7.391667335093565
$ws.Range("E4").Value = 12.31263209739727
$ws.Range("F4").Value = 39.32163919041295
$ws.Range("G4").Value = 46.13114002680524
$ws.Range("H4").Value = 18.52087232557314
$ws.Range("J4").Value = 9.908628130634932
$ws.Range("K4").Value = 16.43267218940929
$ws.Range("M4").Value = 18.43534354674623
$ws.Range("N4").Value = 19.31182189253912
# Row 5
$ws.Range("C5").Value = 4.610702188463699
$ws.Range("D5").Value = 7.390982365483244
$ws.Range("E5").Value = 12.31805203745495
$ws.Range("F5").Value = 39.32498574557949
$ws.Range("G5").Value = 46.10981455375608
$ws.Range("H5").Value = 18.53274818944116
$ws.Range("J5").Value = 9.915423808906652
$ws.Range("K5").Value = 16.34218073438851
$ws.Range("M5").Value = 18.4054034384413
$ws.Range("N5").Value = 19.33126676645111
# Row 6
$ws.Range("C6").Value = 4.603845335833752
$ws.Range("D6").Value = 7.390877819752781
$ws.Range("E6").Value = 12.31897689000885
$ws.Range("F6").Value = 39.32570770357341
$ws.Range("G6").Value = 46.10651380468994
$ws.Range("H6").Value = 18.53478163986943
$ws.Range("J6").Value = 9.916568519812607
$ws.Range("K6").Value = 16.32716476049266
$ws.Range("M6").Value = 18.40048427678995
$ws.Range("N6").Value = 19.33452548499307
# Row 7
$ws.Range("C7").Value = 4.651367271286677
$ws.Range("D7").Value = 7.391657481607441
$ws.Range("E7").Value = 12.3127035250666
$ws.Range("F7").Value = 39.32167317357829
$ws.Range("G7").Value = 46.13083631510892
$ws.Range("H7").Value = 18.52102836458327
$ws.Range("J7").Value = 9.908718687442303
$ws.Range("K7").Value = 16.43145120196255
$ws.Range("M7").Value = 18.43493626905627
$ws.Range("N7").Value = 19.31208212885617
# Row 8
$ws.Range("C8").Value = 4.858761721474189
$ws.Range("D8").Value = 7.396428153806302
$ws.Range("E8").Value = 12.28863525702856
$ws.Range("F8").Value = 39.32807025203237
$ws.Range("G8").Value = 46.27305831583207
$ws.Range("H8").Value = 18.46929572113552
$ws.Range("J8").Value = 9.876445189543352
$ws.Range("K8").Value = 16.89217611403102
$ws.Range("M8").Value = 18.59440068944719
$ws.Range("N8").Value = 19.21737979255067
# Row 9
$ws.Range("C9").Value = 5.252606146126168
$ws.Range("D9").Value = 7.410505256881812
$ws.Range("E9").Value = 12.25441695865691
$ws.Range("F9").Value = 39.42776115322286
$ws.Range("G9").Value = 46.67656056958275
$ws.Range("H9").Value = 18.40024179116513
$ws.Range("J9").Value = 9.821639256876415
$ws.Range("K9").Value = 17.7889966738244
$ws.Range("M9").Value = 18.93204628351199
$ws.Range("N9").Value = 19.04712660155426
# Row 10
$ws.Range("C10").Value = 5.530567642026769
$ws.Range("D10").Value = 7.423672305096388
$ws.Range("E10").Value = 12.23717689605417
$ws.Range("F10").Value = 39.55411399369996
$ws.Range("G10").Value = 47.04779519710304
$ws.Range("H10").Value = 18.36942255080007
$ws.Range("J10").Value = 9.786524299521112
$ws.Range("K10").Value = 18.43659043447173
$ws.Range("M10").Value = 19.19361282699887
$ws.Range("N10").Value = 18.93134529168981
# Row 11
$ws.Range("C11").Value = 5.65390508228547
$ws.Range("D11").Value = 7.430263492957176
$ws.Range("E11").Value = 12.23104586490231
$ws.Range("F11").Value = 39.62307845625725
$ws.Range("G11").Value = 47.232555279359
$ws.Range("H11").Value = 18.35975702636797
$ws.Range("J11").Value = 9.771663925531811
$ws.Range("K11").Value = 18.72732662473103
$ws.Range("M11").Value = 19.31515475002739
$ws.Range("N11").Value = 18.88066792866392
# Row 12
$ws.Range("C12").Value = 5.700116779191081
$ws.Range("D12").Value = 7.432844773162235
$ws.Range("E12").Value = 12.22896992654153
$ws.Range("F12").Value = 39.65083720157936
$ws.Range("G12").Value = 47.30476467206592
$ws.Range("H12").Value = 18.35672511860646
$ws.Range("J12").Value = 9.766196520805353
$ws.Range("K12").Value = 18.83675700724378
$ws.Range("M12").Value = 19.36151352250118
$ws.Range("N12").Value = 18.86176227474083
# Row 13
$ws.Range("C13").Value = 5.690186898488985
$ws.Range("D13").Value = 7.432285071377501
$ws.Range("E13").Value = 12.22940609366965
$ws.Range("F13").Value = 39.64478592380594
$ws.Range("G13").Value = 47.28911392742916
$ws.Range("H13").Value = 18.35735012413652
$ws.Range("J13").Value = 9.76736691651829
$ws.Range("K13").Value = 18.8132204790726
$ws.Range("M13").Value = 19.35151504225088
$ws.Range("N13").Value = 18.86582130764478
# Row 14
$ws.Range("C14").Value = 5.657717085243011
$ws.Range("D14").Value = 7.430474152773549
$ws.Range("E14").Value = 12.23087015331432
$ws.Range("F14").Value = 39.62532930405315
$ws.Range("G14").Value = 47.23845124328998
$ws.Range("H14").Value = 18.35949498825132
$ws.Range("J14").Value = 9.771210915631451
$ws.Range("K14").Value = 18.73634346485765
$ws.Range("M14").Value = 19.31896222573618
$ws.Range("N14").Value = 18.87910685177316
# Row 15
$ws.Range("C15").Value = 5.637762824799215
$ws.Range("D15").Value = 7.429375992067275
$ws.Range("E15").Value = 12.23179892414746
$ws.Range("F15").Value = 39.61362529835736
$ws.Range("G15").Value = 47.20770997241492
$ws.Range("H15").Value = 18.36089064839239
$ws.Range("J15").Value = 9.773586292367725
$ws.Range("K15").Value = 18.68916430272132
$ws.Range("M15").Value = 19.29906510466807
$ws.Range("N15").Value = 18.88728166540929
# Row 16
$ws.Range("C16").Value = 5.522440603458745
$ws.Range("D16").Value = 7.423253560817984
$ws.Range("E16").Value = 12.23761198388289
$ws.Range("F16").Value = 39.54983747510109
$ws.Range("G16").Value = 47.03603706245028
$ws.Range("H16").Value = 18.37014201902485
$ws.Range("J16").Value = 9.787517847882574
$ws.Range("K16").Value = 18.41750291912233
$ws.Range("M16").Value = 19.18571829330053
$ws.Range("N16").Value = 18.93469711759538
# Row 17
$ws.Range("C17").Value = 5.450863221115363
$ws.Range("D17").Value = 7.419650886136693
$ws.Range("E17").Value = 12.24161623261135
$ws.Range("F17").Value = 39.51364236806709
$ws.Range("G17").Value = 46.934763061029
$ws.Range("H17").Value = 18.37693424239676
$ws.Range("J17").Value = 9.796349453822891
$ws.Range("K17").Value = 18.24977867599979
$ws.Range("M17").Value = 19.11681396594632
$ws.Range("N17").Value = 18.96429398045816
# Row 18
$ws.Range("C18").Value = 5.409403644316197
$ws.Range("D18").Value = 7.417635403430816
$ws.Range("E18").Value = 12.24408049617853
$ws.Range("F18").Value = 39.49390564750695
$ws.Range("G18").Value = 46.87801064116099
$ws.Range("H18").Value = 18.38125069128079
$ws.Range("J18").Value = 9.801533990765057
$ws.Range("K18").Value = 18.15295120631588
$ws.Range("M18").Value = 19.07742447173818
$ws.Range("N18").Value = 18.98150492364142
# Row 19
$ws.Range("C19").Value = 5.395317810765056
$ws.Range("D19").Value = 7.416962764806131
$ws.Range("E19").Value = 12.2449425373606
$ws.Range("F19").Value = 39.48740913899448
$ws.Range("G19").Value = 46.85905359657215
$ws.Range("H19").Value = 18.38278247230684
$ws.Range("J19").Value = 9.803307398000397
$ws.Range("K19").Value = 18.12010926722729
$ws.Range("M19").Value = 19.06413052806822
$ws.Range("N19").Value = 18.98736452757777
# Row 20
$ws.Range("C20").Value = 5.458513139264046
$ws.Range("D20").Value = 7.420028538952389
$ws.Range("E20").Value = 12.24117330013948
$ws.Range("F20").Value = 39.51738349523789
$ws.Range("G20").Value = 46.94538911392199
$ws.Range("H20").Value = 18.37616877614392
$ws.Range("J20").Value = 9.795398466467683
$ws.Range("K20").Value = 18.26767097036199
$ws.Range("M20").Value = 19.12412409950758
$ws.Range("N20").Value = 18.96112393936909
# Row 21
$ws.Range("C21").Value = 5.667267995417059
$ws.Range("D21").Value = 7.431003756660042
$ws.Range("E21").Value = 12.23043345713338
$ws.Range("F21").Value = 39.63099966045367
$ws.Range("G21").Value = 47.25327151640501
$ws.Range("H21").Value = 18.35884792495664
$ws.Range("J21").Value = 9.770077502231924
$ws.Range("K21").Value = 18.75894300268262
$ws.Range("M21").Value = 19.32851499785129
$ws.Range("N21").Value = 18.8751968505353
# Row 22
$ws.Range("C22").Value = 5.800806771521789
$ws.Range("D22").Value = 7.438673567756491
$ws.Range("E22").Value = 12.22484655913107
$ws.Range("F22").Value = 39.71482797737878
$ws.Range("G22").Value = 47.46755309355008
$ws.Range("H22").Value = 18.35119023065669
$ws.Range("J22").Value = 9.754460680977582
$ws.Range("K22").Value = 19.07609479878182
$ws.Range("M22").Value = 19.46402410872518
$ws.Range("N22").Value = 18.82069756758911
# Row 23
$ws.Range("C23").Value = 5.729813516448883
$ws.Range("D23").Value = 7.434534969232932
$ws.Range("E23").Value = 12.22769748039933
$ws.Range("F23").Value = 39.66921455064681
$ws.Range("G23").Value = 47.3520059548576
$ws.Range("H23").Value = 18.35494155892561
$ws.Range("J23").Value = 9.762710481040889
$ws.Range("K23").Value = 18.9072184847166
$ws.Range("M23").Value = 19.39153521802222
$ws.Range("N23").Value = 18.8496336188174
# Row 24
$ws.Range("C24").Value = 5.45505556951313
$ws.Range("D24").Value = 7.419857628408703
$ws.Range("E24").Value = 12.2413730447602
$ws.Range("F24").Value = 39.51568879082757
$ws.Range("G24").Value = 46.94058048936524
$ws.Range("H24").Value = 18.37651356176006
$ws.Range("J24").Value = 9.795828074138416
$ws.Range("K24").Value = 18.25958310327826
$ws.Range("M24").Value = 19.12081848619906
$ws.Range("N24").Value = 18.96255650641541
# Row 25
$ws.Range("C25").Value = 5.147826432432174
$ws.Range("D25").Value = 7.406196941370721
$ws.Range("E25").Value = 12.26228547495228
$ws.Range("F25").Value = 39.39145551689777
$ws.Range("G25").Value = 46.55416342183273
$ws.Range("H25").Value = 18.41543654892262
$ws.Range("J25").Value = 9.835559723269949
$ws.Range("K25").Value = 17.54786497166119
$ws.Range("M25").Value = 18.83820634874784
$ws.Range("N25").Value = 19.09154215609399
